$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old demo rows/columns (A1:L4) completely so stale styles/content
# don't linger, then rebuild the new "GuoTai_N" entrust import layout
# (columns A:O, 4 rows incl. header).
$ws.Rows("1:4").Delete()

# --- Header row (row 1): A1:N1 are formulas that just literal-ize the
# column caption text; O1 is a plain "交易类别" label. ---
$ws.Range("A1").Formula = '="委托日期"'
$ws.Range("B1").Formula = '="委托时间"'
$ws.Range("C1").Formula = '="证券代码"'
$ws.Range("D1").Formula = '="证券名称"'
$ws.Range("E1").Formula = '="买卖标志"'
$ws.Range("F1").Formula = '="委托价格"'
$ws.Range("G1").Formula = '="委托数量"'
$ws.Range("H1").Formula = '="委托编号"'
$ws.Range("I1").Formula = '="成交数量"'
$ws.Range("J1").Formula = '="撤单数量"'
$ws.Range("K1").Formula = '="状态说明"'
$ws.Range("L1").Formula = '="撤单标志"'
$ws.Range("M1").Formula = '="股东代码"'
$ws.Range("N1").Formula = '="操作日期"'
$ws.Range("O1").Value = "交易类别"

# --- Data rows 2-4. Columns A, C, D, E, L, M, N repeat the same formula
# across the three rows, so assigning across the whole range in one shot
# reproduces Excel's shared-formula grouping. ---
$ws.Range("A2:A4").Formula = '="20170301"'
$ws.Range("C2:C4").Formula = '="002798"'
$ws.Range("D2:D4").Formula = '="帝王洁具"'
$ws.Range("E2:E4").Formula = '="证券卖出"'
$ws.Range("L2:L4").Formula = '="正常"'
$ws.Range("M2:M4").Formula = '="0208635819"'
$ws.Range("N2:N4").Formula = '="20170301"'

# Row 2
$ws.Range("B2").Formula = '="14:34:37"'
$ws.Range("F2").Value = 58.35
$ws.Range("G2").Value = 4300
$ws.Range("H2").Formula = '="270613"'
$ws.Range("I2").Value = 4300
$ws.Range("J2").Value = 0
$ws.Range("K2").Formula = '="已成"'
$ws.Range("O2").Value = "日内"

# Row 3
$ws.Range("B3").Formula = '="14:35:49"'
$ws.Range("F3").Value = 58.18
$ws.Range("G3").Value = 2000
$ws.Range("H3").Formula = '="271679"'
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 0
$ws.Range("K3").Formula = '="已成"'
$ws.Range("O3").Value = "波段"

# Row 4
$ws.Range("B4").Formula = '="14:37:59"'
$ws.Range("F4").Value = 58.58
$ws.Range("G4").Value = 8000
$ws.Range("H4").Formula = '="273461"'
$ws.Range("I4").Value = 2800
$ws.Range("J4").Value = 5200
$ws.Range("K4").Formula = '="部撤"'
$ws.Range("O4").Value = "目标"

# Match the author's final cursor position (just past the new data block).
[void]$ws.Range("O5").Select()

Write-Output "applied GuoTai_N entrust template edit"
